# ---------------------------------------------------------------------------
# Adds an "Instructions" worksheet (after the existing "Sites" sheet) to the
# MassWateR Sites Template workbook, matching the published template update.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the new worksheet right after "Sites" -------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Instructions"

# --- column widths ----------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 29.451862
$ws2.Columns.Item(2).ColumnWidth = 95.592448
$ws2.Columns.Item(3).ColumnWidth = 20.592448
$ws2.Columns.Item(4).ColumnWidth = 20.736979
$ws2.Columns.Item(5).ColumnWidth = 24.166667

# --- intro notes (rows 1-2) -------------------------------------------------
$ws2.Range("A1").Value2 = "The Sites tab must be formatted exactly like the Sites template, with all of the following fields."
$ws2.Range("A2").Value2 = "The Sites tab must be the first tab in this workbook."
$introRange = $ws2.Range("A1:A2")
$introRange.Font.Bold = $true
$introRange.Font.Color = 12611584

$ws2.Range("C1").Value2 = "Template updated 5/19/23"
$ws2.Range("C1").Font.Color = 255

# --- header row (row 4) ------------------------------------------------------
$ws2.Range("A4").Value2 = "Field"
$ws2.Range("B4").Value2 = "Instructions"
$ws2.Range("C4").Value2 = "Example"
$ws2.Range("D4").Value2 = "Available Values"
$ws2.Range("E4").Value2 = "Required?"

$headerRow = $ws2.Range("A4:E4")
$headerRow.HorizontalAlignment = -4108
$headerRow.Rows.RowHeight = 15.75
$headerRow.Borders.LineStyle = 1
$headerRow.Borders.Weight = 2
$headerRow.Borders.Item(9).Weight = -4138

# --- data rows 5-9 ------------------------------------------------------------
$ws2.Range("A5").Value2 = "Monitoring Location ID"
$ws2.Range("B5").Value2 = "Location ID that is used in your Results file.  Must match exactly."
$ws2.Range("C5").Value2 = "ABT-010"
$ws2.Range("D5").Value2 = "any"
$ws2.Range("E5").Value2 = "Required"

$ws2.Range("A6").Value2 = "Monitoring Location Name"
$ws2.Range("B6").Value2 = "Name of monitoring location."
$ws2.Range("C6").Value2 = "477 Lowell Rd, Concord"
$ws2.Range("D6").Value2 = "any"
$ws2.Range("E6").Value2 = "Required for WQX"

$ws2.Range("A7").Value2 = "Monitoring Location Latitude "
$ws2.Range("B7").Value2 = "Latitude of monitoring location in decimal form.  At least 5 decimals."
$ws2.Range("C7").Value2 = 42.47037
$ws2.Range("D7").Value2 = "any"
$ws2.Range("E7").Value2 = "Required for mapping"

$ws2.Range("A8").Value2 = "Monitoring Location Longitude"
$ws2.Range("B8").Value2 = "Longitude of monitoring location in decimal form.  At least 5 decimals."
$ws2.Range("C8").Value2 = -71.362579
$ws2.Range("D8").Value2 = "any"
$ws2.Range("E8").Value2 = "Required for mapping"

$ws2.Range("A9").Value2 = "Location Group"
$ws2.Range("B9").Value2 = "An optional free-form grouping attribute.  This will allow you to summarize locations by group in the graphing and mapping analysis functions."
$ws2.Range("C9").Value2 = "Lower Assabet"
$ws2.Range("D9").Value2 = "any"
$ws2.Range("E9").Value2 = "Optional"
$ws2.Rows.Item(9).RowHeight = 30

# formatting for the data block A5:E9
$dataBlock = $ws2.Range("A5:E9")
$dataBlock.VerticalAlignment = -4160
$dataBlock.Borders.LineStyle = 1
$dataBlock.Borders.Weight = 2
$ws2.Range("A5").Borders.Item(8).LineStyle = -4142

$ws2.Range("B5:B9").WrapText = $true
$ws2.Range("C5:E9").HorizontalAlignment = -4108
$ws2.Range("D5:E9").Font.Italic = $true

# --- trailing bold empty cell (row 11) --------------------------------------
$ws2.Range("B11").Font.Bold = $true

# --- freeze panes / view state for the Instructions sheet -------------------
$ws2.Activate()
$ws2.Range("B5").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("C2").Select()

# --- restore Sites as the active sheet/selection, per the source workbook ---
$ws1.Activate()
$ws1.Range("E2").Select()
